$wb = $excel.ActiveWorkbook

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 6812797
$ws.Range("I88").Value = 1540
$ws.Range("J88").Value = 10786030
$ws.Range("K88").Value = 1540
$ws.Range("L88").Value = 10786030
$ws.Range("M88").Value = -1134
$ws.Range("N88").Value = -10786842

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 6812797
$ws.Range("I91").Value = 1540
$ws.Range("J91").Value = 10786030
$ws.Range("K91").Value = 1540
$ws.Range("L91").Value = 10786030
$ws.Range("M91").Value = -136
$ws.Range("N91").Value = -10788838

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2565.3845
$ws.Range("I137").Value = 1516.6666
$ws.Range("J137").Value = 2880
$ws.Range("K137").Value = 4549.9998
$ws.Range("L137").Value = 8640
$ws.Range("M137").Value = -1999.9998
$ws.Range("N137").Value = -13740

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4449.54
$ws.Range("I138").Value = 2324.8
$ws.Range("J138").Value = 5360.143
$ws.Range("K138").Value = 6974.400000000001
$ws.Range("L138").Value = 16080.429
$ws.Range("M138").Value = -1834.400000000001
$ws.Range("N138").Value = -26360.429

# ARM row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("N4").Value = 0

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1127334.4
$ws.Range("I32").Value = 12438.389
$ws.Range("J32").Value = 14506087
$ws.Range("K32").Value = 12438.389
$ws.Range("L32").Value = 14506087
$ws.Range("M32").Value = -12151.389
$ws.Range("N32").Value = -14506661

# ARM row 37
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 6303.5
$ws.Range("I37").Value = 500
$ws.Range("J37").Value = 8238
$ws.Range("K37").Value = 500
$ws.Range("L37").Value = 8238
$ws.Range("M37").Value = -227
$ws.Range("N37").Value = -8784

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7912.5
$ws.Range("I61").Value = 10477.75
$ws.Range("J61").Value = 2782
$ws.Range("K61").Value = 10477.75
$ws.Range("L61").Value = 2782
$ws.Range("M61").Value = -10265.75
$ws.Range("N61").Value = -3206

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1980.963
$ws.Range("I74").Value = 2058.88
$ws.Range("J74").Value = 1007
$ws.Range("K74").Value = 2058.88
$ws.Range("L74").Value = 1007
$ws.Range("M74").Value = -1184.88
$ws.Range("N74").Value = -2755

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1980.963
$ws.Range("I77").Value = 2058.88
$ws.Range("J77").Value = 1007
$ws.Range("K77").Value = 10294.4
$ws.Range("L77").Value = 5035
$ws.Range("M77").Value = -5926.400000000001
$ws.Range("N77").Value = -13771

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 29122.264
$ws.Range("I132").Value = 2620.1667
$ws.Range("J132").Value = 128505.125
$ws.Range("K132").Value = 7860.500100000001
$ws.Range("L132").Value = 385515.375
$ws.Range("M132").Value = -5330.500100000001
$ws.Range("N132").Value = -390575.375

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7912.5
$ws.Range("I136").Value = 10477.75
$ws.Range("J136").Value = 2782
$ws.Range("K136").Value = 31433.25
$ws.Range("L136").Value = 8346
$ws.Range("M136").Value = -28883.25
$ws.Range("N136").Value = -13446

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1909.1
$ws.Range("I94").Value = 1795.5714
$ws.Range("J94").Value = 2174
$ws.Range("K94").Value = 1795.5714
$ws.Range("L94").Value = 2174
$ws.Range("M94").Value = -1344.5714
$ws.Range("N94").Value = -3076

# BSM row 138
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 41170
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 41170
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 41170
$ws.Range("N138").Value = -51450

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 927.8316
$ws.Range("I113").Value = 467.85715
$ws.Range("J113").Value = 1007.3333
$ws.Range("K113").Value = 1403.57145
$ws.Range("L113").Value = 3021.9999
$ws.Range("M113").Value = 766.4285500000001
$ws.Range("N113").Value = -7361.9999

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 114074.3
$ws.Range("I80").Value = 187790.5
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 187790.5
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -186792.5
$ws.Range("N80").Value = -5496

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 114074.3
$ws.Range("I83").Value = 187790.5
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 938952.5
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -933960.5
$ws.Range("N83").Value = -27484

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2577
$ws.Range("I102").Value = 900
$ws.Range("J102").Value = 3695
$ws.Range("K102").Value = 900
$ws.Range("L102").Value = 3695
$ws.Range("M102").Value = 722
$ws.Range("N102").Value = -6939

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3242.5715
$ws.Range("I132").Value = 2939.6
$ws.Range("J132").Value = 3518
$ws.Range("K132").Value = 8818.799999999999
$ws.Range("L132").Value = 10554
$ws.Range("M132").Value = -6288.799999999999
$ws.Range("N132").Value = -15614

# GSM row 134
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 36250
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 36250
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 108750
$ws.Range("N134").Value = -113820

# GSM row 136
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 25130.4
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 25130.4
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 75391.20000000001
$ws.Range("N136").Value = -80491.20000000001

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 10988
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 15311.429
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 15311.429
$ws.Range("M46").Value = -712
$ws.Range("N46").Value = -15687.429

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2711.4062
$ws.Range("I61").Value = 2263.0908
$ws.Range("J61").Value = 3697.7
$ws.Range("K61").Value = 2263.0908
$ws.Range("L61").Value = 3697.7
$ws.Range("M61").Value = -2061.0908
$ws.Range("N61").Value = -4101.7

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2711.4062
$ws.Range("I113").Value = 2263.0908
$ws.Range("J113").Value = 3697.7
$ws.Range("K113").Value = 2263.0908
$ws.Range("L113").Value = 3697.7
$ws.Range("M113").Value = -93.09079999999994
$ws.Range("N113").Value = -8037.7

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2711.8794
$ws.Range("I136").Value = 1744.9714
$ws.Range("J136").Value = 4183.2607
$ws.Range("K136").Value = 5234.914199999999
$ws.Range("L136").Value = 12549.7821
$ws.Range("M136").Value = -2684.914199999999
$ws.Range("N136").Value = -17649.7821

# LTW row 137
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5914.5713
$ws.Range("I62").Value = 6222.222
$ws.Range("J62").Value = 5360.8
$ws.Range("K62").Value = 6222.222
$ws.Range("L62").Value = 5360.8
$ws.Range("M62").Value = -5598.222
$ws.Range("N62").Value = -6608.8

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 5914.5713
$ws.Range("I65").Value = 6222.222
$ws.Range("J65").Value = 5360.8
$ws.Range("K65").Value = 31111.11
$ws.Range("L65").Value = 26804
$ws.Range("M65").Value = -27991.11
$ws.Range("N65").Value = -33044

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1494.9412
$ws.Range("I126").Value = 892.38464
$ws.Range("J126").Value = 3453.25
$ws.Range("K126").Value = 2677.15392
$ws.Range("L126").Value = 10359.75
$ws.Range("M126").Value = -207.1539199999997
$ws.Range("N126").Value = -15299.75

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1193.6389
$ws.Range("I136").Value = 806.7692
$ws.Range("J136").Value = 2199.5
$ws.Range("K136").Value = 2420.3076
$ws.Range("L136").Value = 6598.5
$ws.Range("M136").Value = 129.6923999999999
$ws.Range("N136").Value = -11698.5

# WVR row 137
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 70000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 70000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 70000
$ws.Range("N137").Value = -80200

# WVR row 141
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 69137.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 69137.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 69137.5
$ws.Range("N141").Value = -79497.5
